# Updates the "cryptos" worksheet with freshly scraped price / volume data.
# All writes go through a small helper that forces text-typed cells (the
# sheet stores Price/Volume as inline strings, never numbers) while leaving
# the cell's style untouched (no numFmt/style index should be introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row => column letter => new value
$changes = @{
    2  = @{ D = "72.543.54"; E = "  +0.33%  " }
    3  = @{ D = "2.643.73";  E = "  -1.17%  " }
    4  = @{ E = "  +0.01%  " }
    5  = @{ D = "586.17";    E = "  -2.32%  " }
    6  = @{ D = "175.83";    E = "  -0.48%  " }
    8  = @{ D = "0.520";     E = "  -0.58%  " }
    9  = @{ B = "LidoStakedEther"; C = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D = "2.644.44"; E = "  -1.13%  " }
    10 = @{ B = "Dogecoin";        C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge";     D = "0.172";    E = "  +0.94%  " }
    11 = @{ E = "  +1.49%  " }
    12 = @{ E = "  +1.74%  " }
    13 = @{ E = "  -1.62%  " }
    14 = @{ D = "3.129.69";  E = "  -1.24%  " }
    15 = @{ D = "0.0000186"; E = "  +0.47%  " }
    16 = @{ D = "72.358.42"; E = "  +0.22%  " }
    17 = @{ D = "25.83";     E = "  -1.87%  " }
    18 = @{ D = "2.645.54";  E = "  -1.58%  " }
    19 = @{ D = "12.07";     E = "  +0.28%  " }
    20 = @{ D = "375.88";    E = "  +1.46%  " }
    21 = @{ D = "7.86";      E = "  -2.02%  " }
    22 = @{ E = "  -1.33%  " }
    23 = @{ E = "  -0.60%  " }
    24 = @{ D = "71.53";     E = "  -0.61%  " }
    25 = @{ E = "  -0.08%  " }
    26 = @{ E = "  -2.06%  " }
    27 = @{ D = "9.52";      E = "  -3.09%  " }
    28 = @{ D = "2.781.75";  E = "  -1.61%  " }
    29 = @{ D = "0.998";     E = "  -0.22%  " }
    30 = @{ E = "  +1.19%  " }
    31 = @{ E = "  -1.03%  " }
    32 = @{ D = "492.19";    E = "  -3.70%  " }
    33 = @{ D = "1.32";      E = "  +2.05%  " }
    34 = @{ D = "1.80";      E = "  -1.01%  " }
    35 = @{ D = "1.00";      E = "  +0.05%  " }
    36 = @{ D = "161.72";    E = "  -1.11%  " }
    37 = @{ D = "0.116";     E = "  +8.39%  " }
    38 = @{ D = "19.20";     E = "  -1.79%  " }
    39 = @{ E = "  -1.12%  " }
    40 = @{ E = "  -0.96%  " }
    41 = @{ E = "  -0.06%  " }
    42 = @{ E = "  -4.99%  " }
    43 = @{ E = "  +0.31%  " }
    44 = @{ E = "  -2.31%  " }
    45 = @{ E = "  -1.89%  " }
    46 = @{ D = "39.09";     E = "  -0.47%  " }
    47 = @{ D = "150.74";    E = "  -1.66%  " }
    48 = @{ E = "  -1.32%  " }
    49 = @{ E = "  -2.49%  " }
    50 = @{ E = "  -2.55%  " }
    51 = @{ D = "0.609";     E = "  +1.05%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        Set-TextValue $ws.Range($addr) $cols[$col]
    }
}
